# Auto-generated edit script for as-completingthesquare.docx
$d = $word.ActiveDocument

# --- 1) Collapse the multi-run Title/Author/Abstract paragraphs into single runs ---

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Find.Execute("Answers: Completing the square", $true, $false, $false, $false, $false, $true, 1, $false, "Answers: Completing the square", 2) | Out-Null

$authorPara = $d.Paragraphs.Item(2)
$authorPara.Range.Find.Execute("Tom Coleman", $true, $false, $false, $false, $false, $true, 1, $false, "Tom Coleman", 2) | Out-Null

$abstractPara = $d.Paragraphs.Item(4)
$abstractPara.Range.Find.Execute('Answers to questions relating to the guide on completing the square.', $true, $false, $false, $false, $false, $true, 1, $false, 'Answers to questions relating to the guide on completing the square.', 2) | Out-Null

# --- 2) Reorder m:sepChr ahead of m:endChr inside every m:dPr (36 equations) ---

$d.OMaths.Item(1).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>15</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>14</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(4).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>y</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>6</m:t></m:r><m:r><m:t>y</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>8</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>y</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(7).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>8</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>20</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(10).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>m</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>26</m:t></m:r><m:r><m:t>m</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>25</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>m</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>13</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>144</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(13).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>n</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>6</m:t></m:r><m:r><m:t>n</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>50</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>m</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>41</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(16).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>144</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>143</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(19).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>h</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:t>h</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>h</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>4</m:t></m:r></m:den></m:f></m:oMath>') | Out-Null
$d.OMaths.Item(22).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>13</m:t></m:r></m:num><m:den><m:r><m:t>4</m:t></m:r></m:den></m:f></m:oMath>') | Out-Null
$d.OMaths.Item(25).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>13</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>43</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>13</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>4</m:t></m:r></m:den></m:f></m:oMath>') | Out-Null
$d.OMaths.Item(28).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>y</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>8</m:t></m:r><m:r><m:t>y</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>16</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>y</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:oMath>') | Out-Null
$d.OMaths.Item(31).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>13</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>9</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>13</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>133</m:t></m:r></m:num><m:den><m:r><m:t>4</m:t></m:r></m:den></m:f></m:oMath>') | Out-Null
$d.OMaths.Item(34).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>m</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:t>m</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>33</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>m</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>143</m:t></m:r></m:num><m:den><m:r><m:t>4</m:t></m:r></m:den></m:f></m:oMath>') | Out-Null
$d.OMaths.Item(37).Range.InsertXML('<m:oMath><m:r><m:t>2</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>12</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>14</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>2</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(41).Range.InsertXML('<m:oMath><m:r><m:t>5</m:t></m:r><m:sSup><m:e><m:r><m:t>y</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>10</m:t></m:r><m:r><m:t>y</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>5</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(45).Range.InsertXML('<m:oMath><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>32</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>68</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(47).Range.InsertXML('<m:oMath><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>8</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(48).Range.InsertXML('<m:oMath><m:r><m:t>2</m:t></m:r><m:sSup><m:e><m:r><m:t>m</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:t>m</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>2</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>m</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:oMath>') | Out-Null
$d.OMaths.Item(52).Range.InsertXML('<m:oMath><m:r><m:t>3</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>5</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>3</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>14</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f></m:oMath>') | Out-Null
$d.OMaths.Item(56).Range.InsertXML('<m:oMath><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>1</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:oMath>') | Out-Null
$d.OMaths.Item(60).Range.InsertXML('<m:oMath><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:oMath>') | Out-Null
$d.OMaths.Item(61).Range.InsertXML('<m:oMath><m:r><m:t>2</m:t></m:r><m:sSup><m:e><m:r><m:t>h</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:t>h</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>1</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>2</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>h</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>4</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>8</m:t></m:r></m:den></m:f></m:oMath>') | Out-Null
$d.OMaths.Item(65).Range.InsertXML('<m:oMath><m:r><m:t>3</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>5</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>3</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>5</m:t></m:r></m:num><m:den><m:r><m:t>6</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>36</m:t></m:r></m:den></m:f></m:oMath>') | Out-Null
$d.OMaths.Item(69).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>y</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>6</m:t></m:r><m:r><m:t>y</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>8</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>y</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(70).Range.InsertXML('<m:oMath><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>y</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>0</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(75).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>m</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>26</m:t></m:r><m:r><m:t>m</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>25</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>m</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>13</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>144</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(76).Range.InsertXML('<m:oMath><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>y</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>144</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>0</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(81).Range.InsertXML('<m:oMath><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>8</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>20</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(82).Range.InsertXML('<m:oMath><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>±</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:t>i</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(83).Range.InsertXML('<m:oMath><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>0</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(88).Range.InsertXML('<m:oMath><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>1</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:oMath>') | Out-Null
$d.OMaths.Item(89).Range.InsertXML('<m:oMath><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>0</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(92).Range.InsertXML('<m:oMath><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>32</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>68</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(93).Range.InsertXML('<m:oMath><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>±</m:t></m:r><m:r><m:t>i</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(94).Range.InsertXML('<m:oMath><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>0</m:t></m:r></m:oMath>') | Out-Null
$d.OMaths.Item(99).Range.InsertXML('<m:oMath><m:r><m:t>3</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>5</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>3</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>5</m:t></m:r></m:num><m:den><m:r><m:t>6</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>36</m:t></m:r></m:den></m:f></m:oMath>') | Out-Null
$d.OMaths.Item(100).Range.InsertXML('<m:oMath><m:r><m:t>3</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>5</m:t></m:r></m:num><m:den><m:r><m:t>6</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>36</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>0</m:t></m:r></m:oMath>') | Out-Null

Write-Host 'edit complete'
